$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("JIRA_Details")

# Append 5 new rows (37-41) to the JIRA_Details lookup table, carrying the
# same B/C (India/China JIRA id) values forward, mirroring the existing
# formatting of the rows directly above them.
$ws.Range("A36:C36").Copy($ws.Range("A37:C37")) | Out-Null
$ws.Range("A36:C36").Copy($ws.Range("A38:C38")) | Out-Null
$ws.Range("A33:C33").Copy($ws.Range("A39:C39")) | Out-Null
$ws.Range("A33:C33").Copy($ws.Range("A40:C40")) | Out-Null
$ws.Range("A33:C33").Copy($ws.Range("A41:C41")) | Out-Null

$ws.Cells.Item(37,1).Value = "BlockCompanyBrand"
$ws.Cells.Item(38,1).Value = "BlockCompanyClient"
$ws.Cells.Item(39,1).Value = "BlockCompanyProduct"
$ws.Cells.Item(40,1).Value = "BlockGlobalProduct"
$ws.Cells.Item(41,1).Value = "BlockCompanyVendor"

# Make JIRA_Details the active sheet/tab and move the selection onto the
# newly-added data, matching the saved view state.
$ws.Activate()
$ws.Range("A38").Select() | Out-Null
